$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Clear B10's hyperlink-only style BEFORE inserting the new row, so the ---
# --- newly inserted row (which inherits formatting from the row above) picks up ---
# --- plain style 5 instead of the Hyperlink style that currently sits on B10. ---
$ws.Range("B10").Style = "Normal"

# Insert a new row at position 11 (pushes the existing blank row 11 and every
# row below it down by one; the new row inherits formatting from row 10).
$ws.Rows(11).Insert()

# Row 9 content: was ImgFileName / ExtName / screenshot description.
# Becomes the ReceiverMailID row (previously row 10).
$ws.Range("A9").Value = "ReceiverMailID"
$ws.Range("B9").Value = "ganeshsoley@gmail.com"
$ws.Range("C9").Value = "This email ID will receive the mail with snapshot attached."

# Row 10 content: was ReceiverMailID / email / description.
# Becomes the new MailCredential row.
$ws.Range("A10").Value = "MailCredential"
$ws.Range("B10").Value = "GaneshSoley_GMail"
$ws.Range("C10").Value = "Robot will use this value to access mail Credential details from Orchestrator and use them to send mail."

# Row 11 content (brand-new row): MailBodyTemplate row.
$ws.Range("A11").Value = "MailBodyTemplate"
$ws.Range("B11").Value = "Data\Input\MailBody.txt"
$ws.Range("C11").Value = "This template is used in mail to send to the stakeholders."

# --- Hyperlinks ---
# This host's Hyperlinks.Delete() drops the whole sheet collection (and its
# relationships) rather than scoping to the target range/item, so remove all
# of them once and re-add the ones that should still exist: the URL link on
# B6 (unchanged) and the mailto link, which now lives on B9 instead of B10.
$ws.Range("B6").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B6"), "http://www.rpasamples.com/")
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:ganeshsoley@gmail.com")
# Re-apply the named Hyperlink style after Add() so the cell reuses the
# workbook's existing Hyperlink style slot instead of minting a new one.
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("B9").Style = "Hyperlink"

# B10 is plain text now (no hyperlink, no special style).
$ws.Range("B10").Style = "Normal"

# Update the active selection to match the edited workbook.
$ws.Range("C11").Select()
